# Weekly data refresh: a new week's record is inserted at the top of the
# data block (row 140), pushing the existing rows 140-159 down to 141-160.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 140, shifting rows 140:159 down to 141:160.
$ws.Rows.Item(140).Insert()

# Populate the new row 140 with the latest week's record.
$ws.Cells.Item(140, 1).Value = 4
$ws.Cells.Item(140, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(140, 3).Value = "Los Lagos"
$ws.Cells.Item(140, 4).Value = 44491
$ws.Cells.Item(140, 5).Value = 10
$ws.Cells.Item(140, 6).Value = "Fruta"
$ws.Cells.Item(140, 7).Value = 100102
$ws.Cells.Item(140, 8).Value = "Cítricos"
$ws.Cells.Item(140, 9).Value = 100102006
$ws.Cells.Item(140, 10).Value = "Pomelo"
$ws.Cells.Item(140, 11).Value = "Start Ruby"
$ws.Cells.Item(140, 12).Value = "Primera"
$ws.Cells.Item(140, 13).Value = 200
$ws.Cells.Item(140, 14).Value = 11000
$ws.Cells.Item(140, 15).Value = 12000
$ws.Cells.Item(140, 16).Value = 11500
$ws.Cells.Item(140, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(140, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(140, 19).Value = 821
$ws.Cells.Item(140, 20).Value = 14
